$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("O1").Value = "F1 train"

# Column O value updates (rows 2-10, 12-16)
$ws.Range("O2").Value = 0.9411764705882353
$ws.Range("O3").Value = 0.9565217391304348
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 0.7301587301587301
$ws.Range("O6").Value = 0.7368421052631579
$ws.Range("O8").Value = 1
$ws.Range("O9").Value = 0.9859154929577465
$ws.Range("O10").Value = 0.9090909090909091
$ws.Range("O12").Value = 0.8333333333333334
$ws.Range("O13").Value = 0.8857142857142857
$ws.Range("O14").Value = 1
$ws.Range("O15").Value = 1
$ws.Range("O16").Value = 0.6590909090909091

# Row 11 (MLP, 10% feature set) full update
$ws.Range("C11").Value = "{'activation': 'tanh', 'alpha': 0.0001, 'hidden_layer_sizes': (64, 32), 'learning_rate': 'constant'}"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = 6
$ws.Range("I11").Value = 0.55
$ws.Range("J11").Value = 0.4
$ws.Range("K11").Value = 0.3333333333333333
$ws.Range("L11").Value = 0.5
$ws.Range("M11").Value = 0.7272727272727273
$ws.Range("N11").Value = 0.3333333333333333
$ws.Range("O11").Value = 0.9117647058823529
